# Apply the Mateus_Profits numeric cell updates (per the upstream xlsx diff).
# Columns H..N are plain cached values (no formulas in the source file), so
# each changed cell is written directly; cells removed by the diff are cleared
# (ClearContents drops the <c> element entirely, matching the XML diff) and
# cells newly introduced by the diff are created by assigning .Value.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2798.5715
$ws.Range("I32").Value = 3072.5
$ws.Range("J32").Value = 2433.3333
$ws.Range("K32").Value = 3072.5
$ws.Range("L32").Value = 2433.3333
$ws.Range("M32").Value = -2746.5
$ws.Range("N32").Value = -3085.3333
$ws.Range("H43").Value = 3336.5
$ws.Range("I43").Value = 2750
$ws.Range("J43").Value = 3727.5
$ws.Range("K43").Value = 2750
$ws.Range("L43").Value = 3727.5
$ws.Range("M43").Value = -2681
$ws.Range("N43").Value = -3865.5
$ws.Range("H86").Value = 5368.25
$ws.Range("I86").Value = 2931.3333
$ws.Range("K86").Value = 2931.3333
$ws.Range("M86").Value = -1808.3333
$ws.Range("H89").Value = 5368.25
$ws.Range("I89").Value = 2931.3333
$ws.Range("K89").Value = 14656.6665
$ws.Range("M89").Value = -9040.666499999999
$ws.Range("H106").Value = 7511
$ws.Range("I106").Value = 7511
$ws.Range("K106").Value = 7511
$ws.Range("M106").Value = -6880
$ws.Range("H107").Value = 771.125
$ws.Range("I107").Value = 896.8461
$ws.Range("J107").Value = 226.33333
$ws.Range("K107").Value = 896.8461
$ws.Range("L107").Value = 226.33333
$ws.Range("M107").Value = 1023.1539
$ws.Range("N107").Value = -4066.33333
$ws.Range("H109").Value = 39999.5
$ws.Range("I109").Value = 39999
$ws.Range("K109").Value = 39999
$ws.Range("M109").Value = -38612
$ws.Range("H116").Value = 3934
$ws.Range("I116").Value = 3351.3333
$ws.Range("K116").Value = 3351.3333
$ws.Range("M116").Value = 90.66670000000022
$ws.Range("H132").Value = 1637.7916
$ws.Range("I132").Value = 1563.7894
$ws.Range("J132").Value = 1919
$ws.Range("K132").Value = 4691.3682
$ws.Range("L132").Value = 5757
$ws.Range("M132").Value = -2161.3682
$ws.Range("N132").Value = -10817
$ws.Range("H137").Value = 1734
$ws.Range("I137").Value = 1879.6
$ws.Range("K137").Value = 5638.799999999999
$ws.Range("M137").Value = -3088.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 15034.5
$ws.Range("I39").Value = 20069
$ws.Range("K39").Value = 20069
$ws.Range("M39").Value = -19549
$ws.Range("H74").Value = 6661.52
$ws.Range("J74").Value = 9255.111000000001
$ws.Range("L74").Value = 9255.111000000001
$ws.Range("N74").Value = -11003.111
$ws.Range("H77").Value = 6661.52
$ws.Range("J77").Value = 9255.111000000001
$ws.Range("L77").Value = 46275.55500000001
$ws.Range("N77").Value = -55011.55500000001
$ws.Range("H110").Value = 6346.8696
$ws.Range("I110").Value = 4512.1816
$ws.Range("J110").Value = 8028.6665
$ws.Range("K110").Value = 4512.1816
$ws.Range("L110").Value = 8028.6665
$ws.Range("M110").Value = -2467.1816
$ws.Range("N110").Value = -12118.6665
$ws.Range("H132").Value = 4406
$ws.Range("I132").Value = 3732.25
$ws.Range("K132").Value = 11196.75
$ws.Range("M132").Value = -8666.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H94").Value = 2069.75
$ws.Range("I94").Value = 1759.3549
$ws.Range("J94").Value = 3138.889
$ws.Range("K94").Value = 1759.3549
$ws.Range("L94").Value = 3138.889
$ws.Range("M94").Value = -1308.3549
$ws.Range("N94").Value = -4040.889
$ws.Range("H105").Value = 3656.1177
$ws.Range("I105").Value = 3147.4614
$ws.Range("K105").Value = 3147.4614
$ws.Range("M105").Value = -1400.4614
$ws.Range("H118").Value = 40000
$ws.Range("J118").Value = 40000
$ws.Range("L118").Value = 40000
$ws.Range("N118").Value = -43314
$ws.Range("H134").Value = 4765.65
$ws.Range("I134").Value = 4621.737
$ws.Range("J134").Value = 7500
$ws.Range("K134").Value = 13865.211
$ws.Range("L134").Value = 22500
$ws.Range("M134").Value = -11330.211
$ws.Range("N134").Value = -27570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2903.524
$ws.Range("I16").Value = 2377.4666
$ws.Range("K16").Value = 2377.4666
$ws.Range("M16").Value = -2090.4666
$ws.Range("H22").Value = 841.4286
$ws.Range("I22").Value = 298.72726
$ws.Range("J22").Value = 2831.3333
$ws.Range("K22").Value = 298.72726
$ws.Range("L22").Value = 2831.3333
$ws.Range("M22").Value = 51.27274
$ws.Range("N22").Value = -3531.3333
$ws.Range("H94").Value = 2159.4
$ws.Range("J94").Value = 2199.25
$ws.Range("L94").Value = 2199.25
$ws.Range("N94").Value = -3101.25
$ws.Range("H97").Value = 27753
$ws.Range("J97").Value = 27303.6
$ws.Range("L97").Value = 27303.6
$ws.Range("N97").Value = -29285.6
$ws.Range("H107").Value = 1510.2307
$ws.Range("I107").Value = 579.0909
$ws.Range("K107").Value = 579.0909
$ws.Range("M107").Value = 1340.9091
$ws.Range("H113").Value = 2903.524
$ws.Range("I113").Value = 2377.4666
$ws.Range("K113").Value = 2377.4666
$ws.Range("M113").Value = -207.4666000000002
$ws.Range("H122").Value = 2964
$ws.Range("I122").Value = 3114.75
$ws.Range("J122").Value = 2481.6
$ws.Range("K122").Value = 9344.25
$ws.Range("L122").Value = 7444.799999999999
$ws.Range("M122").Value = -6894.25
$ws.Range("N122").Value = -12344.8
$ws.Range("H132").Value = 5304.4
$ws.Range("I132").Value = 3307.3333
$ws.Range("J132").Value = 8300
$ws.Range("K132").Value = 9921.999899999999
$ws.Range("L132").Value = 24900
$ws.Range("M132").Value = -7391.999899999999
$ws.Range("N132").Value = -29960

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 29
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4523.1763
$ws.Range("I80").Value = 2919.2
$ws.Range("J80").Value = 5191.5
$ws.Range("K80").Value = 2919.2
$ws.Range("L80").Value = 5191.5
$ws.Range("M80").Value = -1921.2
$ws.Range("N80").Value = -7187.5
$ws.Range("H83").Value = 4523.1763
$ws.Range("I83").Value = 2919.2
$ws.Range("J83").Value = 5191.5
$ws.Range("K83").Value = 14596
$ws.Range("L83").Value = 25957.5
$ws.Range("M83").Value = -9604
$ws.Range("N83").Value = -35941.5
$ws.Range("H113").Value = 445991
$ws.Range("J113").Value = 649.5
$ws.Range("L113").Value = 649.5
$ws.Range("N113").Value = -4989.5
$ws.Range("H132").Value = 7462.25
$ws.Range("I132").Value = 7462.25
$ws.Range("K132").Value = 22386.75
$ws.Range("M132").Value = -19856.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2232.1765
$ws.Range("I22").Value = 2199.25
$ws.Range("J22").Value = 2242.3076
$ws.Range("K22").Value = 2199.25
$ws.Range("L22").Value = 2242.3076
$ws.Range("M22").Value = -1904.25
$ws.Range("N22").Value = -2832.3076
$ws.Range("H27").Value = 2232.1765
$ws.Range("I27").Value = 2199.25
$ws.Range("J27").Value = 2242.3076
$ws.Range("K27").Value = 2199.25
$ws.Range("L27").Value = 2242.3076
$ws.Range("M27").Value = -2092.25
$ws.Range("N27").Value = -2456.3076
$ws.Range("H46").Value = 5786.871
$ws.Range("I46").Value = 6682.522
$ws.Range("J46").Value = 3211.875
$ws.Range("K46").Value = 6682.522
$ws.Range("L46").Value = 3211.875
$ws.Range("M46").Value = -6494.522
$ws.Range("N46").Value = -3587.875
$ws.Range("H110").Value = 80643.5
$ws.Range("J110").Value = 80643.5
$ws.Range("L110").Value = 80643.5
$ws.Range("N110").Value = -88823.5
$ws.Range("H122").Value = 4745.7144
$ws.Range("I122").Value = 3346
$ws.Range("K122").Value = 10038
$ws.Range("M122").Value = -7588
$ws.Range("H132").Value = 8924.5
$ws.Range("I132").Value = 10194.357
$ws.Range("J132").Value = 4480
$ws.Range("K132").Value = 30583.071
$ws.Range("L132").Value = 13440
$ws.Range("M132").Value = -28053.071
$ws.Range("N132").Value = -18500
$ws.Range("H136").Value = 2250
$ws.Range("I136").Value = 2200
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 6600
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -4050
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 14000
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 18500
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 18500
$ws.Range("M51").Value = -4490
$ws.Range("N51").Value = -19520
$ws.Range("H113").Value = 637.4231
$ws.Range("I113").Value = 600.5
$ws.Range("J113").Value = 720.5
$ws.Range("K113").Value = 1801.5
$ws.Range("L113").Value = 2161.5
$ws.Range("M113").Value = 368.5
$ws.Range("N113").Value = -6501.5
$ws.Range("H126").Value = 3980.353
$ws.Range("I126").Value = 3316.5
$ws.Range("J126").Value = 4928.7144
$ws.Range("K126").Value = 9949.5
$ws.Range("L126").Value = 14786.1432
$ws.Range("M126").Value = -7479.5
$ws.Range("N126").Value = -19726.1432
$ws.Range("H132").Value = 4495.516
$ws.Range("I132").Value = 3768.25
$ws.Range("K132").Value = 11304.75
$ws.Range("M132").Value = -8774.75
$ws.Range("H136").Value = 6645.483
$ws.Range("I136").Value = 4736.1
$ws.Range("K136").Value = 14208.3
$ws.Range("M136").Value = -11658.3

Write-Output "Applied 245 cell updates and 5 clears across 8 sheets"
